$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 - title shape: "Projekt és célja"
#   -> split into 3 runs: "Projekt" (teal highlight), " és " (plain),
#      "célja" (red highlight)
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$titleShape = $slide1.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange

# Sanity: original text is "Projekt és célja" (7 + 4 + 5 = 16 chars)
$run1 = $titleRange.Characters(1, 7)    # "Projekt"
$run2 = $titleRange.Characters(8, 4)    # " és "
$run3 = $titleRange.Characters(12, 5)   # "célja"

# highlight colour 008080 (teal) - VBA RGB(r,g,b) = r + g*256 + b*65536
$run1.Font.Highlight.RGB = 8421376
# highlight colour FF0000 (red)
$run3.Font.Highlight.RGB = 255

# ---------------------------------------------------------------------------
# Slide 2 - content placeholder: remove the placeholder text "kép"
#   -> empty paragraph, no bullet, marL=0 / indent=0
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$picShape = $slide2.Shapes.Item(2)
$picRange = $picShape.TextFrame.TextRange

# Remove the indentation / left margin that comes with the bullet level
$ruler = $picShape.TextFrame.Ruler
$level1 = $ruler.Levels.Item(1)
$level1.LeftMargin = 0
$level1.FirstMargin = 0

# Turn off the bullet for this paragraph
$picRange.ParagraphFormat.Bullet.Type = 0

# Clear the run text itself
$picRange.Text = ""
